$wb = $excel.ActiveWorkbook

# Update the "strain" sheet: row 2 label changes from
# "none_HG105_none_mCh" to "none_R0_none_mCh"
$ws = $wb.Worksheets.Item("strain")
$ws.Range("A2:H2").Value = "none_R0_none_mCh"

# Make "strain" the active/selected sheet (it was "media" before)
$ws.Activate()
$ws.Range("F9").Select()
